# regen sval data to filter save games
#
# The per-stat columns (TB, d2S, K, IP) only take on a small set of
# distinct "s-val" numbers that are looked up per pitching line.
# Regenerating them against the save-game-filtered dataset remaps each
# distinct raw value to a refreshed value; "sum" (column G) is just
# TB + d2S + K + IP recomputed from the refreshed numbers (Win, column F,
# is untouched by the regen).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Map-Tb($v) {
    if ($v -eq 3.642121602190766) { return 3.286832544864788 }
    elseif ($v -eq 0.6717081370404667) { return 0.1190320826869504 }
    elseif ($v -eq 1.539653030276356) { return 0.6606524410359556 }
    elseif ($v -eq 1.016799631930417) { return 0.2917716402565462 }
    elseif ($v -eq 2.349824325366695) { return 1.455362044514542 }
    elseif ($v -eq 0.4405202784308186) { return 0.04271373187048222 }
    elseif ($v -eq 0.2846480661742272) { return 0.01293466051926884 }
    else { return $v }
}

function Map-D2s($v) {
    if ($v -eq 1.307360055219643) { return 1.655778082260271 }
    elseif ($v -eq 0.5246911461486936) { return 0.306821227259698 }
    elseif ($v -eq 0.1990121036566324) { return 0.04071648406533734 }
    else { return $v }
}

function Map-K($v) {
    if ($v -eq 0.4893031109320483) { return 0.7527432677738641 }
    elseif ($v -eq 0.2181074402325811) { return 0.1494219747398047 }
    elseif ($v -eq 1.050893113504925) { return 3.537761648806719 }
    else { return $v }
}

function Map-Ip($v) {
    if ($v -eq 0.345296217628593) { return 0.4942365360607697 }
    else { return $v }
}

for ($r = 2; $r -le 49; $r++) {
    $oldTb  = $ws.Range("B$r").Value2
    $oldD2s = $ws.Range("C$r").Value2
    $oldK   = $ws.Range("D$r").Value2
    $oldIp  = $ws.Range("E$r").Value2

    $newTb  = Map-Tb $oldTb
    $newD2s = Map-D2s $oldD2s
    $newK   = Map-K $oldK
    $newIp  = Map-Ip $oldIp

    $ws.Range("B$r").Value = $newTb
    $ws.Range("C$r").Value = $newD2s
    $ws.Range("D$r").Value = $newK
    $ws.Range("E$r").Value = $newIp
    $ws.Range("G$r").Value = $newTb + $newD2s + $newK + $newIp
}
